# Automatically handle population types for framework plots
# Adds a new "Plots" worksheet (after "Cascades") describing how SIR / UDT
# population-type quantities are aggregated for framework plots, and
# documents the display-name convention via a cell comment.

$wb = $excel.ActiveWorkbook

# --- Add the new "Plots" worksheet as the last tab ----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "Plots"

# --- Populate the table (column order matters for shared-string layout) -
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Quantities"

$ws.Range("A2").Value = "SIR total"
$ws.Range("A3").Value = "UDT total"
$ws.Range("A4").Value = "SIR dict"
$ws.Range("A5").Value = "SIR function"

$ws.Range("B2").Value = "ch_all"
$ws.Range("B3").Value = "all_people"
$ws.Range("B4").Value = "{'alive':['sus','inf','rec']}"
$ws.Range("B5").Value = "{'alive':'sus+inf+rec'}"

# --- Formatting: bold header row, left-aligned code-name-ish columns ----
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("B1").Font.Bold = $true

$ws.Range("A2:A3").HorizontalAlignment = -4131
$ws.Range("B3").HorizontalAlignment = -4131

$ws.Columns.Item(1).ColumnWidth = 13.877604166666666
$ws.Columns.Item(2).ColumnWidth = 20.307291666666664

# --- Documentation comment on the header cell ----------------------------
$commentText = "This column is for the 'display name' of a compartment within a`npopulation cascade, a state that an entity can exist in that is`ndistinct from all other states.`nExamples may include 'Susceptible', 'Infected Stage 1', 'Recovered',`netc.`nIf entities in the network involve two 'orthogonal' descriptors,`ncompartments should combine the status of each state in the title,`ne.g. 'High Income Earner + Year 12 Education', to make sure that each`nentity in a cascade is only ever in one state at a time.`nIt is possible to bundle independent states as analytical features of`ninterest elsewhere in the framework file.`nNote: A display name is a representative label that users interface`nwith (e.g. in databooks and plots).`nIt should be in title or sentence case."
$ws.Range("A1").AddComment($commentText) | Out-Null

# leave the cursor where the author last clicked while reviewing the sheet
$ws.Range("E11").Select()

# --- Tidy up stale selections left on a few other sheets -----------------
$wsCompartments = $wb.Worksheets.Item("Compartments")
$wsCompartments.Activate()
$wsCompartments.Range("A2:A4").Select()

$wsPopTypes = $wb.Worksheets.Item("Population types")
$wsPopTypes.Activate()
$wsPopTypes.Range("A2").Select()

$wsCharacteristics = $wb.Worksheets.Item("Characteristics")
$wsCharacteristics.Activate()
$wsCharacteristics.Range("A11").Select()

# --- Re-activate the new sheet so it is the one shown on open ------------
$ws.Activate()
